# Changes of 24th May 2022
# Update shipment tracking / rate / result columns (P,Q,R) on rows 3-5.
#
# Values that "look like numbers" (tracking numbers, currency strings) must
# still land as literal text (shared-string) cells, not as numeric cells, and
# must keep the sheet's default (unstyled) cell formatting - exactly like the
# surrounding cells already in the sheet. We therefore enter them with a
# leading apostrophe (forces text / quote-prefix) and then reset the cell
# style back to "Normal" so no extra number-format/quote-prefix styling
# sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value2 = "'" + $text
    $range.Style = "Normal"
}

# Row 3
Set-TextValue $ws.Range("P3") "320018621636"

# Row 4
Set-TextValue $ws.Range("P4") "320018621669"
Set-TextValue $ws.Range("Q4") "`$76.67"
Set-TextValue $ws.Range("R4") "FAIL"

# Row 5 (previously blank cells)
Set-TextValue $ws.Range("P5") "320018606875"
Set-TextValue $ws.Range("Q5") "`$43.36"
Set-TextValue $ws.Range("R5") "FAIL"
